$wb = $excel.ActiveWorkbook

$runmanager = $wb.Worksheets.Item("RUNMANAGER")
$data = $wb.Worksheets.Item("DATA")

# --- DATA sheet (sheet2): Row 7 password value changed from @KBTokzan4078 to @KBTokzan2021
$data.Range("E7").Value = "'@KBTokzan2021"

# --- RUNMANAGER sheet (sheet1) ---
# Row 4: execute flag changed from "yes" to "no"
$runmanager.Range("C4").Value = "no"

# New row 7: multiplyOptionsTest
$runmanager.Range("A7").Value = "multiplyOptionsTest"
$runmanager.Range("B7").Value = "To check multiply options test"
$runmanager.Range("C7").Value = "yes"
$runmanager.Range("D7").Value = "'1"
$runmanager.Range("E7").Value = "'1"

# --- DATA sheet (sheet2): new row 11 - multiplyOptionsTest data row
$data.Range("A11").Value = "multiplyOptionsTest"
$data.Range("B11").Value = "yes"
$data.Range("C11").Value = "chrome"
$data.Range("D11").Value = "'"
$data.Range("E11").Value = "'"

# --- Selection / active tab updates ---
$runmanager.Range("B7").Select() | Out-Null
$data.Range("E11").Select() | Out-Null
